$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FPP")

# Populate columns F:I (duplicate block 1) and K:N (duplicate block 2) for rows 9-23
$ws.Range("F9").Value = $ws.Range("A9").Value2
$ws.Range("G9").Value = $ws.Range("B9").Value2
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 15000
$ws.Range("K9").Value = $ws.Range("A9").Value2
$ws.Range("L9").Value = $ws.Range("B9").Value2
$ws.Range("M9").Value = 2
$ws.Range("N9").Value = 15000
$ws.Range("N9").NumberFormat = $ws.Range("D9").NumberFormat

$ws.Range("F10").Value = $ws.Range("A10").Value2
$ws.Range("G10").Value = $ws.Range("B10").Value2
$ws.Range("H10").Value = 9
$ws.Range("I10").Value = 405000
$ws.Range("K10").Value = $ws.Range("A10").Value2
$ws.Range("L10").Value = $ws.Range("B10").Value2
$ws.Range("M10").Value = 3
$ws.Range("N10").Value = 135000
$ws.Range("N10").NumberFormat = $ws.Range("D10").NumberFormat

$ws.Range("F11").Value = $ws.Range("A11").Value2
$ws.Range("G11").Value = $ws.Range("B11").Value2
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 30000
$ws.Range("K11").Value = $ws.Range("A11").Value2
$ws.Range("L11").Value = $ws.Range("B11").Value2
$ws.Range("M11").Value = 4
$ws.Range("N11").Value = 30000
$ws.Range("N11").NumberFormat = $ws.Range("D11").NumberFormat

$ws.Range("F12").Value = $ws.Range("A12").Value2
$ws.Range("G12").Value = $ws.Range("B12").Value2
$ws.Range("H12").Value = 9
$ws.Range("I12").Value = 405000
$ws.Range("K12").Value = $ws.Range("A12").Value2
$ws.Range("L12").Value = $ws.Range("B12").Value2
$ws.Range("M12").Value = 3
$ws.Range("N12").Value = 135000
$ws.Range("N12").NumberFormat = $ws.Range("D12").NumberFormat

$ws.Range("F13").Value = $ws.Range("A13").Value2
$ws.Range("G13").Value = $ws.Range("B13").Value2
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 30000
$ws.Range("K13").Value = $ws.Range("A13").Value2
$ws.Range("L13").Value = $ws.Range("B13").Value2
$ws.Range("M13").Value = 4
$ws.Range("N13").Value = 30000
$ws.Range("N13").NumberFormat = $ws.Range("D13").NumberFormat

$ws.Range("F14").Value = $ws.Range("A14").Value2
$ws.Range("G14").Value = $ws.Range("B14").Value2
$ws.Range("H14").Value = 6
$ws.Range("I14").Value = 270000
$ws.Range("K14").Value = $ws.Range("A14").Value2
$ws.Range("L14").Value = $ws.Range("B14").Value2
$ws.Range("M14").Value = 2
$ws.Range("N14").Value = 90000
$ws.Range("N14").NumberFormat = $ws.Range("D14").NumberFormat

$ws.Range("F15").Value = $ws.Range("A15").Value2
$ws.Range("G15").Value = $ws.Range("B15").Value2
$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 15000
$ws.Range("K15").Value = $ws.Range("A15").Value2
$ws.Range("L15").Value = $ws.Range("B15").Value2
$ws.Range("M15").Value = 2
$ws.Range("N15").Value = 15000
$ws.Range("N15").NumberFormat = $ws.Range("D15").NumberFormat

$ws.Range("F16").Value = $ws.Range("A16").Value2
$ws.Range("G16").Value = $ws.Range("B16").Value2
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 135000
$ws.Range("K16").Value = $ws.Range("A16").Value2
$ws.Range("L16").Value = $ws.Range("B16").Value2
$ws.Range("M16").Value = 1
$ws.Range("N16").Value = 45000
$ws.Range("N16").NumberFormat = $ws.Range("D16").NumberFormat

$ws.Range("F17").Value = $ws.Range("A17").Value2
$ws.Range("G17").Value = $ws.Range("B17").Value2
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 7500
$ws.Range("K17").Value = $ws.Range("A17").Value2
$ws.Range("L17").Value = $ws.Range("B17").Value2
$ws.Range("M17").Value = 1
$ws.Range("N17").Value = 7500
$ws.Range("N17").NumberFormat = $ws.Range("D17").NumberFormat

$ws.Range("F18").Value = $ws.Range("A18").Value2
$ws.Range("G18").Value = $ws.Range("B18").Value2
$ws.Range("H18").Value = 3
$ws.Range("I18").Value = 135000
$ws.Range("K18").Value = $ws.Range("A18").Value2
$ws.Range("L18").Value = $ws.Range("B18").Value2
$ws.Range("M18").Value = 1
$ws.Range("N18").Value = 45000
$ws.Range("N18").NumberFormat = $ws.Range("D18").NumberFormat

$ws.Range("F19").Value = $ws.Range("A19").Value2
$ws.Range("G19").Value = $ws.Range("B19").Value2
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 7500
$ws.Range("K19").Value = $ws.Range("A19").Value2
$ws.Range("L19").Value = $ws.Range("B19").Value2
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = 7500
$ws.Range("N19").NumberFormat = $ws.Range("D19").NumberFormat

$ws.Range("F20").Value = $ws.Range("A20").Value2
$ws.Range("G20").Value = $ws.Range("B20").Value2
$ws.Range("H20").Value = 6
$ws.Range("I20").Value = 270000
$ws.Range("K20").Value = $ws.Range("A20").Value2
$ws.Range("L20").Value = $ws.Range("B20").Value2
$ws.Range("M20").Value = 2
$ws.Range("N20").Value = 90000
$ws.Range("N20").NumberFormat = $ws.Range("D20").NumberFormat

$ws.Range("F21").Value = $ws.Range("A21").Value2
$ws.Range("G21").Value = $ws.Range("B21").Value2
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = 15000
$ws.Range("K21").Value = $ws.Range("A21").Value2
$ws.Range("L21").Value = $ws.Range("B21").Value2
$ws.Range("M21").Value = 2
$ws.Range("N21").Value = 15000
$ws.Range("N21").NumberFormat = $ws.Range("D21").NumberFormat

$ws.Range("F22").Value = $ws.Range("A22").Value2
$ws.Range("G22").Value = $ws.Range("B22").Value2
$ws.Range("H22").Value = 3
$ws.Range("I22").Value = 135000
$ws.Range("K22").Value = $ws.Range("A22").Value2
$ws.Range("L22").Value = $ws.Range("B22").Value2
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 45000
$ws.Range("N22").NumberFormat = $ws.Range("D22").NumberFormat

$ws.Range("F23").Value = $ws.Range("A23").Value2
$ws.Range("G23").Value = $ws.Range("B23").Value2
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 22500
$ws.Range("K23").Value = $ws.Range("A23").Value2
$ws.Range("L23").Value = $ws.Range("B23").Value2
$ws.Range("M23").Value = 3
$ws.Range("N23").Value = 22500
$ws.Range("N23").NumberFormat = $ws.Range("D23").NumberFormat


# Row 24 totals
$ws.Range("C24").Formula = "=SUM(C9:C23)"
$ws.Range("I24").Formula = "=SUM(I9:I23)"
$ws.Range("M24:N24").Formula = "=SUM(M9:M23)"
$ws.Range("N24").NumberFormat = $ws.Range("D24").NumberFormat

# Column widths for K, L, N (bestFit columns) - mirror the source columns'
# bestFit widths (L mirrors B, N mirrors D); K gets its own bestFit width.
$ws.Columns.Item(12).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(4).ColumnWidth
$ws.Columns.Item(11).ColumnWidth = 9.65

# Selection state
$r = $ws.Range("M22,M10,M12,M14,M16,M18,M20")
$r.Select()
